$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4526511508623453
$ws.Range("C2").Value = 0.5888436429448646
$ws.Range("D2").Value = 0.5812924534180843
$ws.Range("E2").Value = 0.5801507848226425
$ws.Range("F2").Value = 0.5974208531964308
$ws.Range("G2").Value = 0.5433951162106502
$ws.Range("H2").Value = 0.5916991841078396
$ws.Range("I2").Value = 0.5263655152051964
$ws.Range("J2").Value = 0.5763023614572208
$ws.Range("K2").Value = 0.6063861008393496
$ws.Range("L2").Value = 0.5644507163064624
$ws.Range("M2").Value = 0.5666219507292947
$ws.Range("N2").Value = 0.5774081479548764
$ws.Range("O2").Value = 0.5644507163064624
